$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.167.13"
$ws.Range("E2").Value = "  +0.14%  "
$ws.Range("D3").Value = "3.274.43"
$ws.Range("E3").Value = "  +0.57%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "584.26"
$ws.Range("E5").Value = "  -0.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "185.07"
$ws.Range("E6").Value = "  +2.01%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  +1.15%  "
$ws.Range("E9").Value = "  -4.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.62"
$ws.Range("E10").Value = "  -0.54%  "
$ws.Range("E11").Value = "  -2.77%  "
$ws.Range("D12").Value = "3.842.30"
$ws.Range("E12").Value = "  +0.49%  "
$ws.Range("E13").Value = "  +0.92%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.56"
$ws.Range("E14").Value = "  -2.35%  "
$ws.Range("D15").Value = "68.164.23"
$ws.Range("E15").Value = "  +0.12%  "
$ws.Range("E16").Value = "  -1.49%  "
$ws.Range("D17").Value = "3.295.58"
$ws.Range("E17").Value = "  +1.55%  "
$ws.Range("E18").Value = "  -1.53%  "
$ws.Range("E19").Value = "  -1.15%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "417.41"
$ws.Range("E20").Value = "  +6.14%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.55"
$ws.Range("E21").Value = "  -1.62%  "
$ws.Range("E22").Value = "  -0.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "71.54"
$ws.Range("E23").Value = "  +0.17%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.509"
$ws.Range("E24").Value = "  -1.61%  "
$ws.Range("E25").Value = "  -1.42%  "
$ws.Range("E26").Value = "  -1.03%  "
$ws.Range("E27").Value = "  -1.83%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  +0.01%  "
$ws.Range("E29").Value = "  -1.76%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.74"
$ws.Range("E30").Value = "  -1.22%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.48"
$ws.Range("E31").Value = "  -3.35%  "
$ws.Range("E32").Value = "  -3.08%  "
$ws.Range("E33").Value = "  +0.04%  "
$ws.Range("E34").Value = "  -1.90%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "163.90"
$ws.Range("E35").Value = "  -0.55%  "
$ws.Range("E36").Value = "  -2.62%  "
$ws.Range("E37").Value = "  -1.01%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "27.17"
$ws.Range("E38").Value = "  +2.66%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.800"
$ws.Range("E39").Value = "  -2.65%  "
$ws.Range("E40").Value = "  -2.79%  "
$ws.Range("E41").Value = "  -3.59%  "
$ws.Range("D42").Value = "2.665.86"
$ws.Range("E42").Value = "  +2.60%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "40.92"
$ws.Range("E43").Value = "  -1.17%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0680"
$ws.Range("E44").Value = "  -1.41%  "
$ws.Range("E45").Value = "  -1.25%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "337.48"
$ws.Range("E46").Value = "  -1.62%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "24.55"
$ws.Range("E47").Value = "  -0.63%  "
$ws.Range("E48").Value = "  -2.52%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.30"
$ws.Range("E49").Value = "  -0.19%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.979"
$ws.Range("E50").Value = "  +0.05%  "
$ws.Range("E51").Value = "  -1.64%  "
